$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.616.59"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").Value = "1.665.91"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.89"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4803"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2631"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06160"
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").Value = "1.663.19"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.82"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5934"
$ws.Range("E13").Value = "  -3.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.397"
$ws.Range("E14").Value = "  -3.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.47"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.0000"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "25.610.43"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006767"
$ws.Range("E19").Value = "  +3.00%  "
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "1.876.61"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.445"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.312"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.87"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.09"
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.403"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "104.76"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.694"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.958"
$ws.Range("E30").Value = "  +5.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.672"
$ws.Range("E31").Value = "  +4.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07657"
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9994"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04321"
$ws.Range("E34").Value = "  -5.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.616"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6126"
$ws.Range("E36").Value = "  +6.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9524"
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.610"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8571"
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9999"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01510"
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.885"
$ws.Range("E42").Value = "  +3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.07"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3773"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.713"
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.223"
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05268"
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.52"
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.367"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("E51").Value = "  +0.03%  "
